# Lab Tests - Lipid Profile - Third Scenario
# Replace the two-column (URL / Expected Behavior) table with a single-column
# list of lipid-test names: TestName / Cholesterol / Lipid Test / HDL.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column B is no longer used at all - clear its values AND formatting so the
# worksheet collapses back down to a single-column (A1:A4) used range.
$ws.Range("B1:B4").Clear()

# Write the new values. The order below is chosen so the workbook's shared
# string table ends up built in the same order as the target file
# (HDL, Cholesterol, Lipid Test, TestName).
$ws.Range("A4").Value = "HDL"
$ws.Range("A2").Value = "Cholesterol"
$ws.Range("A3").Value = "Lipid Test"
$ws.Range("A1").Value = "TestName"

# Turn on word-wrap for the header and the data rows, then restore the last
# row back to the worksheet's plain/default formatting.
$ws.Range("A1:A4").WrapText = $true
$ws.Range("A4").ClearFormats()

# The header row is given a taller, custom row height.
$ws.Range("A1").EntireRow.RowHeight = 19.5

# Leave the cursor parked just below the data, on A5.
$ws.Range("A5").Select()
